# Applies the cryptos.xlsx price/volume refresh described by the commit
# "Updated cryptos list on Mon Sep  9 04:41:49 UTC 2024 with GitHub Actions".
# Only column D (Price) and column E (Volume(1h)) text values change; every
# value in this sheet is stored as text (note the padded "  +x.xx%  " strings
# and the dotted price format), so plain numeric-looking prices are forced
# back to text (NumberFormat "@" then restoring the "Normal" style so no
# stray cell formatting is left behind) to avoid Excel re-typing them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "54.684.08"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "2.279.81"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "504.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.529"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "2.298.79"
$ws.Range("E9").Value = "  +0.48%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0968"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("E12").Value = "  +3.39%  "
$ws.Range("E13").Value = "  +4.28%  "
$ws.Range("E14").Value = "  +4.41%  "
$ws.Range("D15").Value = "2.685.96"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").Value = "54.758.36"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "2.303.63"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.35"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "307.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.993"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  +0.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.76%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "171.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.72%  "
$ws.Range("E29").Value = "  +2.12%  "
$ws.Range("E30").Value = "  +0.50%  "
$ws.Range("E31").Value = "  +1.93%  "
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.95"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.15%  "
$ws.Range("E35").Value = "  -0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.912"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.78%  "
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.26%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "126.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "252.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.21%  "
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0902"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.35%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.549"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.375"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.82%  "
$ws.Range("E50").Value = "  +0.18%  "
$ws.Range("E51").Value = "  +0.50%  "
